# Append three new daily rows (212-214) to Sheet1, continuing the existing
# series: row 211 (2025-04-20) is duplicated verbatim - formatting and all
# metric values (B:J) unchanged - into rows 212-214, with column A (the date
# serial) advanced by one day each time (2025-04-21, 22, 23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRange = $ws.Range("A211:J211")

$newRows = @(212, 213, 214)
$dateSerials = @(45768, 45769, 45770)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $newRows[$i]
    $destRange = $ws.Range("A" + $row + ":J" + $row)

    # Copy the whole source row (values + number format/font/border/alignment)
    # straight into the destination row so the new cells share row 211's
    # existing style instead of minting a new one.
    $srcRange.Copy($destRange)

    # Column A is the one value that actually changes: advance the date.
    $ws.Cells.Item($row, 1).Value = $dateSerials[$i]
}
